# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# OFF sheet - update Road ("R") row with Wild Card round stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 304
$wsOff.Range("C3").Value = 210
$wsOff.Range("D3").Value = 58
$wsOff.Range("E3").Value = 23
$wsOff.Range("F3").Value = 4
$wsOff.Range("G3").Value = 4

# DEF sheet - update Road ("R") row with Wild Card round stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 287
$wsDef.Range("C3").Value = 210
$wsDef.Range("D3").Value = 67
$wsDef.Range("E3").Value = 33
$wsDef.Range("F3").Value = 4
$wsDef.Range("G3").Value = 5
